$wb = $excel.ActiveWorkbook

# Sheets using the "Ano YYYY" label pattern for the header row (B1:E1)
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet using the "Intervalo ..." label pattern for the header row (B1:E1)
$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIntervalo.Range("B1").Value = "Intervalo 2015"
$wsIntervalo.Range("C1").Value = "Intervalo 2015-2030"
$wsIntervalo.Range("D1").Value = "Intervalo 2031-2040"
$wsIntervalo.Range("E1").Value = "Intervalo 2041-2050"

# Sheet with only a single year label in the header row (B1)
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano 2015"
